$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 214
$ws.Range("I2").Value = 626
$ws.Range("J2").Value = 2410
$ws.Range("L2").Value = 675
$ws.Range("M2").Value = 39
$ws.Range("N2").Value = 449
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 17
$ws.Range("Q2").Value = 1
$ws.Range("R2").Value = 29
$ws.Range("S2").Value = 259
$ws.Range("T2").Value = 433
$ws.Range("U2").Value = 35
$ws.Range("V2").Value = 3806
$ws.Range("W2").Value = 0
$ws.Range("X2").Value = 3693
$ws.Range("Y2").Value = 7
$ws.Range("Z2").Value = 58
$ws.Range("AA2").Value = 18
